$wb = $excel.ActiveWorkbook

# --- Update the daily conversion text on sheet "Hoja1" ---
$wsHoja1 = $wb.Worksheets.Item("Hoja1")
$cell = $wsHoja1.Range("A1")
$text = $cell.Value()
$text = $text -replace "6195\.74 pesos = 1\.69 = 913\.34 Bs", "6212.77 pesos = 1.69 = 931.07 Bs"
$text = $text -replace "1000 Bs = 1\.7 = 6195\.74 pesos", "1000 Bs = 1.7 = 6212.77 pesos"
$cell.Value = $text

# --- Update the rate values on sheet "tasas" ---
$wsTasas = $wb.Worksheets.Item("tasas")
$wsTasas.Range("O10").Value = 3650
$wsTasas.Range("O12").Value = 550
